$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

# --- Fill in missing description / severity details for the new checks ---
# Row 17: Subroutine signature mismatch (severity stays Error)
$ws.Range("D17").Value = "Subroutine signature mismatch"

# Row 18: Variable shadows outer scope -> severity becomes Warning
$ws.Range("A18").Value = "Warning"
$ws.Range("D18").Value = "Variable shadows outer scope"

# Row 19: Variable may be used before initialization
$ws.Range("D19").Value = "Variable may be used before initialization"

# --- Convert the data range into a proper Excel Table ---
$rng = $ws.Range("A1:F28")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium5"

$ws.Range("D29").Select()
